# Auto-generated Excel COM-interop script
# Applies the "Fix Training Data Issue (#48)" data corrections to Sheet1
# - Corrects ~136 numeric stat values that were off due to a 1-day shift in source data
# - Rewrites the BF (Date) column from "2-13-2013-14" to the correct ISO date "2014-02-13"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric corrections -------------------------------------------------
$numericUpdates = @(
    @{ Cell = "AD2"; Value = 25 },
    @{ Cell = "AH2"; Value = 8 },
    @{ Cell = "BA2"; Value = 26 },
    @{ Cell = "AD3"; Value = 2 },
    @{ Cell = "AF3"; Value = 26 },
    @{ Cell = "D4"; Value = 50 },
    @{ Cell = "F4"; Value = 26 },
    @{ Cell = "G4"; Value = 0.48 },
    @{ Cell = "I4"; Value = 35.4 },
    @{ Cell = "J4"; Value = 78.5 },
    @{ Cell = "K4"; Value = 0.451 },
    @{ Cell = "L4"; Value = 8.1 },
    @{ Cell = "M4"; Value = 22 },
    @{ Cell = "N4"; Value = 0.368 },
    @{ Cell = "P4"; Value = 24.5 },
    @{ Cell = "Q4"; Value = 0.764 },
    @{ Cell = "R4"; Value = 9.300000000000001 },
    @{ Cell = "S4"; Value = 30.1 },
    @{ Cell = "T4"; Value = 39.4 },
    @{ Cell = "U4"; Value = 20.9 },
    @{ Cell = "Y4"; Value = 4.1 },
    @{ Cell = "AB4"; Value = 97.7 },
    @{ Cell = "AC4"; Value = -2 },
    @{ Cell = "AD4"; Value = 30 },
    @{ Cell = "AF4"; Value = 14 },
    @{ Cell = "AG4"; Value = 16 },
    @{ Cell = "AN4"; Value = 11 },
    @{ Cell = "AP4"; Value = 9 },
    @{ Cell = "AU4"; Value = 19 },
    @{ Cell = "BB4"; Value = 21 },
    @{ Cell = "BC4"; Value = 20 },
    @{ Cell = "AP5"; Value = 8 },
    @{ Cell = "AU5"; Value = 20 },
    @{ Cell = "D6"; Value = 51 },
    @{ Cell = "E6"; Value = 26 },
    @{ Cell = "G6"; Value = 0.51 },
    @{ Cell = "J6"; Value = 80.8 },
    @{ Cell = "K6"; Value = 0.424 },
    @{ Cell = "M6"; Value = 17.6 },
    @{ Cell = "N6"; Value = 0.336 },
    @{ Cell = "Q6"; Value = 0.768 },
    @{ Cell = "R6"; Value = 12.1 },
    @{ Cell = "S6"; Value = 32.7 },
    @{ Cell = "U6"; Value = 22.1 },
    @{ Cell = "W6"; Value = 7.4 },
    @{ Cell = "X6"; Value = 5.2 },
    @{ Cell = "Y6"; Value = 6.4 },
    @{ Cell = "Z6"; Value = 19 },
    @{ Cell = "AC6"; Value = -0.4 },
    @{ Cell = "AD6"; Value = 25 },
    @{ Cell = "AK6"; Value = 29 },
    @{ Cell = "AU6"; Value = 12 },
    @{ Cell = "BC6"; Value = 16 },
    @{ Cell = "AK7"; Value = 28 },
    @{ Cell = "AN7"; Value = 16 },
    @{ Cell = "AD8"; Value = 2 },
    @{ Cell = "AD9"; Value = 25 },
    @{ Cell = "AF9"; Value = 16 },
    @{ Cell = "AG9"; Value = 18 },
    @{ Cell = "AD10"; Value = 15 },
    @{ Cell = "AZ12"; Value = 10 },
    @{ Cell = "AD13"; Value = 15 },
    @{ Cell = "D15"; Value = 52 },
    @{ Cell = "F15"; Value = 34 },
    @{ Cell = "G15"; Value = 0.346 },
    @{ Cell = "I15"; Value = 37.1 },
    @{ Cell = "L15"; Value = 9.1 },
    @{ Cell = "M15"; Value = 24.6 },
    @{ Cell = "O15"; Value = 17.1 },
    @{ Cell = "Q15"; Value = 0.758 },
    @{ Cell = "S15"; Value = 32.8 },
    @{ Cell = "U15"; Value = 23.1 },
    @{ Cell = "V15"; Value = 15.1 },
    @{ Cell = "Y15"; Value = 4.5 },
    @{ Cell = "AB15"; Value = 100.4 },
    @{ Cell = "AD15"; Value = 15 },
    @{ Cell = "AI15"; Value = 19 },
    @{ Cell = "AM15"; Value = 5 },
    @{ Cell = "AQ15"; Value = 16 },
    @{ Cell = "AS15"; Value = 10 },
    @{ Cell = "AV15"; Value = 17 },
    @{ Cell = "AY15"; Value = 13 },
    @{ Cell = "AZ15"; Value = 11 },
    @{ Cell = "BA15"; Value = 27 },
    @{ Cell = "AD16"; Value = 15 },
    @{ Cell = "AD17"; Value = 25 },
    @{ Cell = "AH17"; Value = 8 },
    @{ Cell = "AM17"; Value = 14 },
    @{ Cell = "AQ17"; Value = 17 },
    @{ Cell = "AV17"; Value = 18 },
    @{ Cell = "AD18"; Value = 15 },
    @{ Cell = "AM18"; Value = 17 },
    @{ Cell = "AG19"; Value = 17 },
    @{ Cell = "AS19"; Value = 9 },
    @{ Cell = "AD20"; Value = 15 },
    @{ Cell = "AD21"; Value = 15 },
    @{ Cell = "AN21"; Value = 12 },
    @{ Cell = "BB21"; Value = 22 },
    @{ Cell = "D22"; Value = 54 },
    @{ Cell = "E22"; Value = 42 },
    @{ Cell = "G22"; Value = 0.778 },
    @{ Cell = "J22"; Value = 82.2 },
    @{ Cell = "K22"; Value = 0.476 },
    @{ Cell = "M22"; Value = 20.4 },
    @{ Cell = "N22"; Value = 0.36 },
    @{ Cell = "O22"; Value = 19.3 },
    @{ Cell = "Q22"; Value = 0.802 },
    @{ Cell = "R22"; Value = 11.1 },
    @{ Cell = "T22"; Value = 45.5 },
    @{ Cell = "W22"; Value = 8.1 },
    @{ Cell = "Z22"; Value = 21.6 },
    @{ Cell = "AA22"; Value = 20.1 },
    @{ Cell = "AB22"; Value = 104.9 },
    @{ Cell = "AC22"; Value = 7.5 },
    @{ Cell = "AD22"; Value = 2 },
    @{ Cell = "AJ22"; Value = 20 },
    @{ Cell = "AM22"; Value = 18 },
    @{ Cell = "AN22"; Value = 15 },
    @{ Cell = "AR22"; Value = 16 },
    @{ Cell = "AD23"; Value = 2 },
    @{ Cell = "AJ23"; Value = 19 },
    @{ Cell = "AD24"; Value = 2 },
    @{ Cell = "AU24"; Value = 11 },
    @{ Cell = "AD25"; Value = 25 },
    @{ Cell = "AV25"; Value = 19 },
    @{ Cell = "AM26"; Value = 4 },
    @{ Cell = "AF27"; Value = 26 },
    @{ Cell = "AG27"; Value = 27 },
    @{ Cell = "AI27"; Value = 18 },
    @{ Cell = "AD29"; Value = 15 },
    @{ Cell = "AY29"; Value = 14 },
    @{ Cell = "AD30"; Value = 15 },
    @{ Cell = "AR30"; Value = 17 },
    @{ Cell = "AD31"; Value = 15 },
    @{ Cell = "AF31"; Value = 16 },
    @{ Cell = "BC31"; Value = 15 }
)

foreach ($u in $numericUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Date column correction ------------------------------------------------
# The original text "2-13-2013-14" is ambiguous/wrong; Excel auto-detects
# plain ISO-like strings as dates, so we force Text format before assigning
# and then clear the transient formatting so the cell keeps its original
# (unstyled) appearance.
$dateCells = @(
    "BF2",
    "BF3",
    "BF4",
    "BF5",
    "BF6",
    "BF7",
    "BF8",
    "BF9",
    "BF10",
    "BF11",
    "BF12",
    "BF13",
    "BF14",
    "BF15",
    "BF16",
    "BF17",
    "BF18",
    "BF19",
    "BF20",
    "BF21",
    "BF22",
    "BF23",
    "BF24",
    "BF25",
    "BF26",
    "BF27",
    "BF28",
    "BF29",
    "BF30",
    "BF31"
)

$correctDate = "2014-02-13"
foreach ($cellRef in $dateCells) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $correctDate
    $rng.ClearFormats()
}

